$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Swap the data in rows 37 and 38 (columns F:V) ---
# Row 37 becomes "Rosario Central vs Talleres Cordoba" (previously row 38's match)
# Row 38 becomes "Argentinos Jrs vs Atl. Tucuman" (previously row 37's match)
$ws.Cells.Item(37, 6).Value = "Rosario Central"
$ws.Cells.Item(37, 7).Value = 2
$ws.Cells.Item(37, 8).Value = "Talleres Cordoba"
$ws.Cells.Item(37, 9).Value = 0
$ws.Cells.Item(37, 10).Value = 2.67
$ws.Cells.Item(37, 11).Value = "29/08/2023 00:12"
$ws.Cells.Item(37, 12).Value = 3.38
$ws.Cells.Item(37, 13).Value = "03/09/2023 21:07"
$ws.Cells.Item(37, 14).Value = 3.13
$ws.Cells.Item(37, 15).Value = "29/08/2023 00:12"
$ws.Cells.Item(37, 16).Value = 3.18
$ws.Cells.Item(37, 17).Value = "03/09/2023 21:13"
$ws.Cells.Item(37, 18).Value = 2.73
$ws.Cells.Item(37, 19).Value = "29/08/2023 00:12"
$ws.Cells.Item(37, 20).Value = 2.36
$ws.Cells.Item(37, 21).Value = "03/09/2023 21:13"
$ws.Cells.Item(37, 22).Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/rosario-central-talleres-cordoba/Q71g1Ok0/"

$ws.Cells.Item(38, 6).Value = "Argentinos Jrs"
$ws.Cells.Item(38, 7).Value = 2
$ws.Cells.Item(38, 8).Value = "Atl. Tucuman"
$ws.Cells.Item(38, 9).Value = 2
$ws.Cells.Item(38, 10).Value = 1.76
$ws.Cells.Item(38, 11).Value = "29/08/2023 01:42"
$ws.Cells.Item(38, 12).Value = 1.75
$ws.Cells.Item(38, 13).Value = "03/09/2023 21:13"
$ws.Cells.Item(38, 14).Value = 3.56
$ws.Cells.Item(38, 15).Value = "29/08/2023 01:42"
$ws.Cells.Item(38, 16).Value = 3.46
$ws.Cells.Item(38, 17).Value = "03/09/2023 21:13"
$ws.Cells.Item(38, 18).Value = 5.34
$ws.Cells.Item(38, 19).Value = "29/08/2023 01:42"
$ws.Cells.Item(38, 20).Value = 5.63
$ws.Cells.Item(38, 21).Value = "03/09/2023 21:13"
$ws.Cells.Item(38, 22).Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/argentinos-jrs-atl-tucuman/hf3o3pLm/"

# --- Step 2: Add 6 new match rows (62-67), copying formatting from row 61 ---
$ws.Range("A61:V61").Copy()
$ws.Range("A62:V67").PasteSpecial(-4122)
$ws.Cells.Item(62, 1).Value = 61
$ws.Cells.Item(62, 2).Value = "argentina"
$ws.Cells.Item(62, 3).Value = "copa-de-la-liga-profesional"
$ws.Cells.Item(62, 4).NumberFormat = "@"
$ws.Cells.Item(62, 4).Value = "2023"
$ws.Cells.Item(62, 4).Style = "Normal"
$ws.Cells.Item(62, 5).Value = 45189.875
$ws.Cells.Item(62, 6).Value = "Barracas Central"
$ws.Cells.Item(62, 7).Value = 1
$ws.Cells.Item(62, 8).Value = "Banfield"
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 10).Value = 2.72
$ws.Cells.Item(62, 11).Value = "16/09/2023 23:12"
$ws.Cells.Item(62, 12).Value = 3.41
$ws.Cells.Item(62, 13).Value = "20/09/2023 20:59"
$ws.Cells.Item(62, 14).Value = 2.92
$ws.Cells.Item(62, 15).Value = "16/09/2023 23:12"
$ws.Cells.Item(62, 16).Value = 2.8
$ws.Cells.Item(62, 17).Value = "20/09/2023 20:52"
$ws.Cells.Item(62, 18).Value = 3
$ws.Cells.Item(62, 19).Value = "16/09/2023 23:12"
$ws.Cells.Item(62, 20).Value = 2.61
$ws.Cells.Item(62, 21).Value = "20/09/2023 20:59"
$ws.Cells.Item(62, 22).Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/barracas-central-banfield/UDN4fDYp/"

$ws.Cells.Item(63, 1).Value = 62
$ws.Cells.Item(63, 2).Value = "argentina"
$ws.Cells.Item(63, 3).Value = "copa-de-la-liga-profesional"
$ws.Cells.Item(63, 4).NumberFormat = "@"
$ws.Cells.Item(63, 4).Value = "2023"
$ws.Cells.Item(63, 4).Style = "Normal"
$ws.Cells.Item(63, 5).Value = 45189.875
$ws.Cells.Item(63, 6).Value = "Defensa y Justicia"
$ws.Cells.Item(63, 7).Value = 2
$ws.Cells.Item(63, 8).Value = "Tigre"
$ws.Cells.Item(63, 9).Value = 0
$ws.Cells.Item(63, 10).Value = 2.49
$ws.Cells.Item(63, 11).Value = "15/09/2023 23:13"
$ws.Cells.Item(63, 12).Value = 2.67
$ws.Cells.Item(63, 13).Value = "20/09/2023 20:58"
$ws.Cells.Item(63, 14).Value = 3.2
$ws.Cells.Item(63, 15).Value = "15/09/2023 23:13"
$ws.Cells.Item(63, 16).Value = 3.03
$ws.Cells.Item(63, 17).Value = "20/09/2023 20:53"
$ws.Cells.Item(63, 18).Value = 3.05
$ws.Cells.Item(63, 19).Value = "15/09/2023 23:13"
$ws.Cells.Item(63, 20).Value = 3.04
$ws.Cells.Item(63, 21).Value = "20/09/2023 20:53"
$ws.Cells.Item(63, 22).Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/defensa-y-justicia-tigre/0bg4FEAA/"

$ws.Cells.Item(64, 1).Value = 63
$ws.Cells.Item(64, 2).Value = "argentina"
$ws.Cells.Item(64, 3).Value = "copa-de-la-liga-profesional"
$ws.Cells.Item(64, 4).NumberFormat = "@"
$ws.Cells.Item(64, 4).Value = "2023"
$ws.Cells.Item(64, 4).Style = "Normal"
$ws.Cells.Item(64, 5).Value = 45189.97916666666
$ws.Cells.Item(64, 6).Value = "Estudiantes L.P."
$ws.Cells.Item(64, 7).Value = 0
$ws.Cells.Item(64, 8).Value = "San Lorenzo"
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 2.04
$ws.Cells.Item(64, 11).Value = "16/09/2023 21:12"
$ws.Cells.Item(64, 12).Value = 2.29
$ws.Cells.Item(64, 13).Value = "20/09/2023 23:29"
$ws.Cells.Item(64, 14).Value = 3.04
$ws.Cells.Item(64, 15).Value = "16/09/2023 21:12"
$ws.Cells.Item(64, 16).Value = 2.85
$ws.Cells.Item(64, 17).Value = "20/09/2023 23:29"
$ws.Cells.Item(64, 18).Value = 4.47
$ws.Cells.Item(64, 19).Value = "16/09/2023 21:12"
$ws.Cells.Item(64, 20).Value = 4.05
$ws.Cells.Item(64, 21).Value = "20/09/2023 23:29"
$ws.Cells.Item(64, 22).Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/estudiantes-l-p-san-lorenzo/IZPUVZfi/"

$ws.Cells.Item(65, 1).Value = 64
$ws.Cells.Item(65, 2).Value = "argentina"
$ws.Cells.Item(65, 3).Value = "copa-de-la-liga-profesional"
$ws.Cells.Item(65, 4).NumberFormat = "@"
$ws.Cells.Item(65, 4).Value = "2023"
$ws.Cells.Item(65, 4).Style = "Normal"
$ws.Cells.Item(65, 5).Value = 45189.97916666666
$ws.Cells.Item(65, 6).Value = "Velez Sarsfield"
$ws.Cells.Item(65, 7).Value = 2
$ws.Cells.Item(65, 8).Value = "Arsenal Sarandi"
$ws.Cells.Item(65, 9).Value = 1
$ws.Cells.Item(65, 10).Value = 1.69
$ws.Cells.Item(65, 11).Value = "17/09/2023 23:43"
$ws.Cells.Item(65, 12).Value = 1.74
$ws.Cells.Item(65, 13).Value = "20/09/2023 23:23"
$ws.Cells.Item(65, 14).Value = 3.41
$ws.Cells.Item(65, 15).Value = "17/09/2023 23:43"
$ws.Cells.Item(65, 16).Value = 3.4
$ws.Cells.Item(65, 17).Value = "20/09/2023 23:23"
$ws.Cells.Item(65, 18).Value = 5.46
$ws.Cells.Item(65, 19).Value = "17/09/2023 23:43"
$ws.Cells.Item(65, 20).Value = 6
$ws.Cells.Item(65, 21).Value = "20/09/2023 23:27"
$ws.Cells.Item(65, 22).Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/velez-sarsfield-arsenal-sarandi/OEeXCLlC/"

$ws.Cells.Item(66, 1).Value = 65
$ws.Cells.Item(66, 2).Value = "argentina"
$ws.Cells.Item(66, 3).Value = "copa-de-la-liga-profesional"
$ws.Cells.Item(66, 4).NumberFormat = "@"
$ws.Cells.Item(66, 4).Value = "2023"
$ws.Cells.Item(66, 4).Style = "Normal"
$ws.Cells.Item(66, 5).Value = 45190.08333333334
$ws.Cells.Item(66, 6).Value = "Instituto"
$ws.Cells.Item(66, 7).Value = 3
$ws.Cells.Item(66, 8).Value = "Colon Santa Fe"
$ws.Cells.Item(66, 9).Value = 1
$ws.Cells.Item(66, 10).Value = 2.04
$ws.Cells.Item(66, 11).Value = "17/09/2023 01:12"
$ws.Cells.Item(66, 12).Value = 2.18
$ws.Cells.Item(66, 13).Value = "21/09/2023 01:56"
$ws.Cells.Item(66, 14).Value = 3.19
$ws.Cells.Item(66, 15).Value = "17/09/2023 01:12"
$ws.Cells.Item(66, 16).Value = 3
$ws.Cells.Item(66, 17).Value = "21/09/2023 01:56"
$ws.Cells.Item(66, 18).Value = 4.18
$ws.Cells.Item(66, 19).Value = "17/09/2023 01:12"
$ws.Cells.Item(66, 20).Value = 4.13
$ws.Cells.Item(66, 21).Value = "21/09/2023 01:56"
$ws.Cells.Item(66, 22).Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/instituto-colon-santa-fe/A90iIGQi/"

$ws.Cells.Item(67, 1).Value = 66
$ws.Cells.Item(67, 2).Value = "argentina"
$ws.Cells.Item(67, 3).Value = "copa-de-la-liga-profesional"
$ws.Cells.Item(67, 4).NumberFormat = "@"
$ws.Cells.Item(67, 4).Value = "2023"
$ws.Cells.Item(67, 4).Style = "Normal"
$ws.Cells.Item(67, 5).Value = 45190.08333333334
$ws.Cells.Item(67, 6).Value = "Racing Club"
$ws.Cells.Item(67, 7).Value = 2
$ws.Cells.Item(67, 8).Value = "Newells Old Boys"
$ws.Cells.Item(67, 9).Value = 1
$ws.Cells.Item(67, 10).Value = 2.62
$ws.Cells.Item(67, 11).Value = "16/09/2023 21:12"
$ws.Cells.Item(67, 12).Value = 2.35
$ws.Cells.Item(67, 13).Value = "21/09/2023 01:55"
$ws.Cells.Item(67, 14).Value = 2.94
$ws.Cells.Item(67, 15).Value = "16/09/2023 21:12"
$ws.Cells.Item(67, 16).Value = 3.13
$ws.Cells.Item(67, 17).Value = "21/09/2023 01:55"
$ws.Cells.Item(67, 18).Value = 3.12
$ws.Cells.Item(67, 19).Value = "16/09/2023 21:12"
$ws.Cells.Item(67, 20).Value = 3.47
$ws.Cells.Item(67, 21).Value = "21/09/2023 01:55"
$ws.Cells.Item(67, 22).Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/racing-club-newells-old-boys/4Wym7jIj/"

Write-Output "Edit complete"
